# Update points for phone 79174466 -> 0.00
# Adds a new row (16) to the customers sheet:
#   A16 = "79174466" (kept as text, matching source data which stores this
#         phone number as a string rather than a number)
#   B16 = "" (empty birthday, like the other rows with unknown birthdays)
#   C16 = 0 (total_points reset to 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text instead of
# auto-converting the numeric-looking string to a number. Resetting the
# style back to "Normal" afterwards drops the quote-prefix formatting that
# Excel applies when using the apostrophe trick, so the new cells stay
# unstyled just like the rest of the data rows.
$ws.Range("A16").Value = "'79174466"
$ws.Range("B16").Value = "'"
$ws.Range("A16:B16").Style = "Normal"

$ws.Range("C16").Value = 0
